$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric keep their exact literal text
# (e.g. trailing zeros like '1.00' or '6.40') by forcing a Text number format
# before assigning the value, matching the source data's plain-text representation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.313.35'
$ws.Range("E2").Value = '  +4.99%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.448.66'
$ws.Range("E3").Value = '  +3.08%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.79'
$ws.Range("E5").Value = '  +4.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.57'
$ws.Range("E6").Value = '  +8.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.517'
$ws.Range("E7").Value = '  +2.54%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  +9.81%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.03'
$ws.Range("E10").Value = '  +4.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0806'
$ws.Range("E11").Value = '  +1.98%  '

$ws.Range("E12").Value = '  -2.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.54'
$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.02'
$ws.Range("E14").Value = '  +3.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.832.51'
$ws.Range("E15").Value = '  +3.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.421.33'
$ws.Range("E16").Value = '  +2.87%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.837'
$ws.Range("E17").Value = '  +3.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.198.22'
$ws.Range("E18").Value = '  +4.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.42'
$ws.Range("E19").Value = '  +3.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.40'
$ws.Range("E20").Value = '  +0.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0928'
$ws.Range("E21").Value = '  +4.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.22'
$ws.Range("E22").Value = '  +1.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.49'
$ws.Range("E23").Value = '  +3.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.29'
$ws.Range("E24").Value = '  +2.58%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.52'
$ws.Range("E25").Value = '  +3.16%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.53'
$ws.Range("E27").Value = '  +4.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("E28").Value = '  -7.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.56'
$ws.Range("E29").Value = '  +2.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.80'
$ws.Range("E30").Value = '  +5.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.51'
$ws.Range("E31").Value = '  +3.11%  '

$ws.Range("E32").Value = '  +15.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.48'
$ws.Range("E33").Value = '  +14.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.25'
$ws.Range("E34").Value = '  +3.57%  '

$ws.Range("E35").Value = '  +0.23%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0766'
$ws.Range("E36").Value = '  +3.41%  '

$ws.Range("E37").Value = '  +4.55%  '

$ws.Range("E38").Value = '  +4.30%  '

$ws.Range("E39").Value = '  +0.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '124.92'
$ws.Range("E40").Value = '  -2.53%  '

$ws.Range("E41").Value = '  +2.57%  '

$ws.Range("E42").Value = '  -2.15%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.45'
$ws.Range("E43").Value = '  +2.55%  '

$ws.Range("E44").Value = '  +4.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.942.37'
$ws.Range("E45").Value = '  +0.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.97'
$ws.Range("E46").Value = '  +6.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.10'
$ws.Range("E47").Value = '  -1.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.27'
$ws.Range("E48").Value = '  +0.36%  '

$ws.Range("E49").Value = '  +17.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.34'
$ws.Range("E50").Value = '  +6.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.94'
$ws.Range("E51").Value = '  +3.93%  '
